# upgrade to v 1.0.16
# Add a "Total" row under the Age column data: a new "Total" label in C12
# and a SUM formula in E12 that totals the Age entries (E7:E11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C12").Value = "Total"
$ws.Range("E12").Formula = "=SUM(E7:E11)"
